# Fix the map dark forest terrain: row 5 (13000002 / 昏暗密林)
#  - add a Quest value (column F) for the row
#  - change the TilePath (column I) from "default" to "darkforest"
#  - move the active selection from F4 to F6

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scene")

$ws.Range("F5").Value = "42000006;2|42000007;1|42000008;2|42000003;2|42000004;1"
$ws.Range("I5").Value = "darkforest"

$ws.Range("F6").Select()
